$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")
$lo = $ws.ListObjects.Item(1)

# Insert a new blank column before "en_variable-label" (column E), shifting
# the existing columns E:P to F:Q, and grow the "Table5" table to include it.
$ws.Range("E1").EntireColumn.Insert()
$lo.Resize($ws.Range("A1:Q5"))

# Name the newly inserted table column "timeval".
$ws.Range("E1").Value2 = "timeval"

# For the "time" HEADING row, the variable-type ("TIME") moves into the new
# timeval column as a boolean flag, and the old variable-type cell is cleared.
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value2 = $true

# Match the column width of the newly inserted column to its left neighbour,
# but without the "best fit" auto-sizing flag (as Excel leaves it after a
# manual column insert).
$ws.Columns(5).ColumnWidth = $ws.Columns(4).ColumnWidth

# Update selection / view state to reflect where the edit was made.
$ws.Range("D4").Select()
